# DOMA-1858: parsing of "Reading submission date" column value either in
# "YYYY-MM-DD" or "YYYY-MM" format.
#
# Update the example values in the "Дата передачи показаний" (Reading
# submission date) column (K) so the sample file demonstrates both
# accepted date formats: the full "YYYY-MM-DD" form and the shorter
# "YYYY-MM" form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = "2021-12"
$ws.Range("K9").Value = "2021-12"
$ws.Range("K10").Value = "2021-12"
